# Apply the cyclic update to rows 4, 5, 8, 9, 11, 13 on the active sheet.
# Each row takes on the player/position/team values that (in the prior
# version of the sheet) belonged to the next row in this sequence, with
# row 4 wrapping around to take what row 13 used to hold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Gradey Dick"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Toronto Raptors"

$ws.Range("A5").Value = "Brandon Ingram"
$ws.Range("B5").Value = "SG,SF,PF"
$ws.Range("C5").Value = "New Orleans Pelicans"

$ws.Range("A8").Value = "Cameron Johnson"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Brooklyn Nets"

$ws.Range("A9").Value = "Anthony Davis"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Los Angeles Lakers"

$ws.Range("A11").Value = "Mason Plumlee"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Phoenix Suns"

$ws.Range("A13").Value = "LaMelo Ball"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Charlotte Hornets"
